$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.647.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.41%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.242.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.67%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.38%  "

$ws.Range("E6").Value = "  +0.92%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.242.21"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.64%  "

$ws.Range("E10").Value = "  +1.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.29%  "

$ws.Range("E12").Value = "  -1.51%  "

$ws.Range("E13").Value = "  +2.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.776.94"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.73%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.690.16"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.245.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.57%  "

$ws.Range("E18").Value = "  -1.32%  "

$ws.Range("E19").Value = "  +1.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "508.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.747"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.92%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.72%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.82%  "

$ws.Range("B26").Value = "Hedera"
$ws.Range("C26").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.174"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +94.37%  "

$ws.Range("E27").Value = "  +0.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.21%  "

$ws.Range("E30").Value = "  -1.97%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.80%  "

$ws.Range("E33").Value = "  +0.67%  "

$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("E35").Value = "  -4.61%  "

$ws.Range("E36").Value = "  -2.84%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0802"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +18.80%  "

$ws.Range("E38").Value = "  +1.02%  "

$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +15.82%  "

$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "492.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.61%  "

$ws.Range("E41").Value = "  +0.44%  "

$ws.Range("E42").Value = "  +1.85%  "

$ws.Range("E43").Value = "  -1.28%  "

$ws.Range("E44").Value = "  -2.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.24%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.938.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.19%  "

$ws.Range("E49").Value = "  +1.79%  "

$ws.Range("E51").Value = "  -0.82%  "
